$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A84:B84").Copy()
$ws.Range("A85:B85").PasteSpecial(-4122)

$ws.Range("A85").Value = "MaxAllowedDeviceExceeded"
$ws.Range("B85").Value = "TestOkur’u sadece kayıtlı olan bilgisayarınızda kullanabilirsiniz. Daha fazla bilgisayarda kullanabilmek için yeni lisans satın almanız gerekmektedir"
